$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, matching the formatting of the
# existing header cells (copy G1's format onto H1, then set its text).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add corresponding value in H2 (plain numeric cell, no special style)
$ws.Range("H2").Value = 1
